# Rename the header row of the products sheet (Лист1 / sheet1).
#
# Old headers: article | name | quantity | size
# New headers: art     | description | qty | pack
#
# Data rows (article numbers, "Eau de Parfum" label, quantity, size) are
# left untouched - only the column header labels change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "qty"
$ws.Range("D1").Value = "pack"
$ws.Range("A1").Value = "art"
